$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Header: rename "Features" (G1) to "Obs" ---
$ws.Range("G1").Value = "Obs"

# --- 2) Insert a new row before row 7 (Random Forest block gains a row) ---
$ws.Rows(7).EntireRow.Insert()
$ws.Range("A6:I6").Copy()
$ws.Range("A7:I7").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G7").Value = "Select from model (logistic regression)"
$ws.Range("H7").Value = 0.99
$ws.Range("I7").Value = 0.98199999999999998

# --- 3) Insert a new row before (current) row 12 (SVM block gains a row) ---
$ws.Rows(12).EntireRow.Insert()
$ws.Range("A13:I13").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G12").Value = "Select from model (logistic regression)"
$ws.Range("H12").Value = 0.99
$ws.Range("I12").Value = 0.98199999999999998

# --- 4) Restore view/selection state ---
$ws.Range("A4").Select()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("I13").Select()
